# TestDataShareSkill.xlsx - add more ShareSkill test rows + tweak sheet views
$wb = $excel.ActiveWorkbook

$shareSkill = $wb.Worksheets.Item("ShareSkill")
$deleteShareSkill = $wb.Worksheets.Item("DeleteShareSkill")

# --- New data rows on ShareSkill ---
$shareSkill.Range("B5").Value = "Design website banners and assist with web visuals"
$shareSkill.Range("C5").Value = "Digital Marketing"
$shareSkill.Range("E5").Value = "Advertising "
$shareSkill.Range("F5").Value = "Hourly basis service"
$shareSkill.Range("G5").Value = "On-site"
$shareSkill.Range("M5").Value = "Skill-Exchange"
$shareSkill.Range("P5").Value = "Active"

$shareSkill.Range("B6").Value = "Design website banners and assist with web visuals"
$shareSkill.Range("C6").Value = "Digital Marketing"
$shareSkill.Range("D6").Value = "Social Media Marketing"
$shareSkill.Range("E6").Value = "Advertising "
$shareSkill.Range("F6").Value = "Hourly basis service"
$shareSkill.Range("G6").Value = "On-site"
$shareSkill.Range("M6").Value = "Skill-Exchange"
$shareSkill.Range("P6").Value = "Active"

$shareSkill.Range("B5:B6").WrapText = $true
$shareSkill.Range("B5:B6").VerticalAlignment = -4160

$shareSkill.Range("A7").Value = "Testtesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttesttestt105"
$shareSkill.Range("B7").Value = " "

$shareSkill.Range("A8").Value = "Testtesttesttestte20"

$shareSkill.Range("B9").Value = "Testing Software and Web ApplicationsTesting 555 Testing Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web Applications403Testing Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web Ap600603"

$shareSkill.Range("B10").Value = "Testing Software and Web ApplicationsTesting 555 Testing Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web ApplicationsTesting Software and Web Applications197"

# --- Column width tweak on ShareSkill (col D) ---
$shareSkill.Columns.Item(4).ColumnWidth = 22.33203125

# --- Sheet view changes ---
$shareSkill.Range("F5").Select()
$ws0 = $excel.ActiveWindow
$ws0.ScrollColumn = 3

$deleteShareSkill.Range("A3").Select()

$shareSkill.Activate()

$excel.ActiveWindow.WindowState = -4143
